$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for columns G:T (rows 2-17)
$rowData = @{
    2 = @(21.139235, 63.417705, 0.1633331201667119, 0.1633331201667119, 3, 1, 8.658142333333332, 25.974427, 0.1344393815966713, 0.1344393815966713, 183.0265054477816, 1647.238549030035, 0.02195840366946755, 0.02195840366946755)
    3 = @(21.139235, 63.417705, 0.1633331201667119, 0.1633331201667119, 3, 1, 2.796453333333333, 8.38936, 0.04342195384682981, 0.04342195384682981, 59.11488417986667, 532.0339576187999, 0.007092243205537673, 0.007092243205537674)
    4 = @(21.139235, 63.417705, 0.1633331201667119, 0.1633331201667119, 3, 1, 4.152187333333333, 12.456562, 0.06447312551305152, 0.06447312551305154, 87.77406380335665, 789.96657423021, 0.01053059675694674, 0.01053059675694675)
    5 = @(21.139235, 63.417705, 0.1633331201667119, 0.1633331201667119, 3, 1, 48.795048, 146.385144, 0.7576655390434474, 0.7576655390434474, 1031.48998650828, 9283.409878574519, 0.12375187653476, 0.12375187653476)
    6 = @(74.39645633333333, 223.189369, 0.5748271090353965, 0.5748271090353966, 3, 1, 8.658142333333332, 25.974427, 0.1344393815966713, 0.1344393815966713, 644.135108029618, 5797.215972266563, 0.07727940106372103, 0.07727940106372104)
    7 = @(74.39645633333333, 223.189369, 0.5748271090353965, 0.5748271090353966, 3, 1, 2.796453333333333, 8.38936, 0.04342195384682981, 0.04342195384682981, 208.0462183015378, 1872.41596471384, 0.02496011619844159, 0.0249601161984416)
    8 = @(74.39645633333333, 223.189369, 0.5748271090353965, 0.5748271090353966, 3, 1, 4.152187333333333, 12.456562, 0.06447312551305152, 0.06447312551305154, 308.9080236321531, 2780.172212689378, 0.03706090034914367, 0.03706090034914369)
    9 = @(74.39645633333333, 223.189369, 0.5748271090353965, 0.5748271090353966, 3, 1, 48.795048, 146.385144, 0.7576655390434474, 0.7576655390434474, 3630.178657814904, 32671.60792033414, 0.4355266914240902, 0.4355266914240902)
    10 = @(4.054539666666667, 12.163619, 0.03132755819197652, 0.03132755819197652, 3, 1, 8.658142333333332, 25.974427, 0.1344393815966713, 0.1344393815966713, 35.10478153014589, 315.943033771313, 0.004211657550263057, 0.004211657550263057)
    11 = @(4.054539666666667, 12.163619, 0.03132755819197652, 0.03132755819197652, 3, 1, 2.796453333333333, 8.38936, 0.04342195384682981, 0.04342195384682981, 11.33833096598222, 102.04497869384, 0.00136030378594588, 0.00136030378594588)
    12 = @(4.054539666666667, 12.163619, 0.03132755819197652, 0.03132755819197652, 3, 1, 4.152187333333333, 12.456562, 0.06447312551305152, 0.06447312551305154, 16.83520824643089, 151.516874217878, 0.002019785591328728, 0.002019785591328728)
    13 = @(4.054539666666667, 12.163619, 0.03132755819197652, 0.03132755819197652, 3, 1, 48.795048, 146.385144, 0.7576655390434474, 0.7576655390434474, 197.841457652904, 1780.573118876136, 0.02373581126443886, 0.02373581126443886)
    14 = @(29.83382566666667, 89.50147699999999, 0.2305122126059151, 0.2305122126059151, 3, 1, 8.658142333333332, 25.974427, 0.1344393815966713, 0.1344393815966713, 258.3055089698532, 2324.749580728679, 0.03098991931321963, 0.03098991931321964)
    15 = @(29.83382566666667, 89.50147699999999, 0.2305122126059151, 0.2305122126059151, 3, 1, 2.796453333333333, 8.38936, 0.04342195384682981, 0.04342195384682981, 83.42890123163556, 750.86011108472, 0.01000929065690466, 0.01000929065690467)
    16 = @(29.83382566666667, 89.50147699999999, 0.2305122126059151, 0.2305122126059151, 3, 1, 4.152187333333333, 12.456562, 0.06447312551305152, 0.06447312551305154, 123.8756330380082, 1114.880697342074, 0.01486184281563238, 0.01486184281563238)
    17 = @(29.83382566666667, 89.50147699999999, 0.2305122126059151, 0.2305122126059151, 3, 1, 48.795048, 146.385144, 0.7576655390434474, 0.7576655390434474, 1455.742955428632, 13101.68659885769, 0.1746511598201584, 0.1746511598201584)
}

foreach ($rnum in $rowData.Keys) {
    $vals = $rowData[$rnum]
    $col = 7  # column G
    foreach ($v in $vals) {
        $ws.Cells.Item($rnum, $col).Value = $v
        $col = $col + 1
    }
}
